# Update the build timestamp embedded in the version string throughout the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---------------------------------------------------------

$aboutSheet.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$aboutSheet.Range("A6").Value = "Recommended Citation:  " + '"' + "Global Energy Monitor, Coal mine boundaries and methane sources for Huashan Coal Mine, China, M1965, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)" + '"'

# --- "Boundaries and methane sources" sheet --------------------------------

for ($row = 2; $row -le 12; $row++) {
    $cell = $dataSheet.Range("S$row")
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
}
